$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.634.19'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '3.563.06'
$ws.Range("E3").Value = '  +0.79%  '

$ws.Range("E4").Value = '  +0.09%  '

$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").Value = "'607.11"
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = '  -0.07%  '

$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").Value = "'145.18"
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = '  +0.78%  '

$ws.Range("D7").Value = '3.562.68'
$ws.Range("E7").Value = '  +0.84%  '

$ws.Range("E8").Value = '  +0.09%  '

$style_D9 = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.498"
$ws.Range("D9").Style = $style_D9
$ws.Range("E9").Value = '  +3.70%  '

$style_D11 = $ws.Range("D11").Style
$ws.Range("D11").Value = "'7.97"
$ws.Range("D11").Style = $style_D11
$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("E12").Value = '  +0.94%  '

$ws.Range("D13").Value = '4.170.70'
$ws.Range("E13").Value = '  +0.89%  '

$ws.Range("E14").Value = '  -0.02%  '

$ws.Range("E15").Value = '  -0.55%  '

$ws.Range("D16").Value = '3.587.47'
$ws.Range("E16").Value = '  +1.52%  '

$ws.Range("D17").Value = '66.714.88'
$ws.Range("E17").Value = '  +0.49%  '

$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").Value = "'0.116"
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = '  +0.36%  '

$style_D19 = $ws.Range("D19").Style
$ws.Range("D19").Value = "'11.52"
$ws.Range("D19").Style = $style_D19
$ws.Range("E19").Value = '  +5.08%  '

$ws.Range("E20").Value = '  +0.00%  '

$style_D21 = $ws.Range("D21").Style
$ws.Range("D21").Value = "'14.93"
$ws.Range("D21").Style = $style_D21
$ws.Range("E21").Value = '  +0.04%  '

$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").Value = "'432.18"
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = '  +1.50%  '

$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").Value = "'0.618"
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = '  +2.79%  '

$style_D24 = $ws.Range("D24").Style
$ws.Range("D24").Value = "'79.80"
$ws.Range("D24").Style = $style_D24
$ws.Range("E24").Value = '  +1.29%  '

$ws.Range("D25").Value = '3.709.33'
$ws.Range("E25").Value = '  +0.89%  '

$ws.Range("E26").Value = '  -0.22%  '

$ws.Range("E27").Value = '  -0.20%  '

$style_D28 = $ws.Range("D28").Style
$ws.Range("D28").Value = "'8.04"
$ws.Range("D28").Style = $style_D28
$ws.Range("E28").Value = '  -1.09%  '

$ws.Range("E29").Value = '  +1.31%  '

$ws.Range("E30").Value = '  -0.28%  '

$ws.Range("E31").Value = '  +0.04%  '

$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").Value = "'1.46"
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = '  -1.60%  '

$ws.Range("D33").Value = '3.560.95'
$ws.Range("E33").Value = '  +0.99%  '

$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").Value = "'25.36"
$ws.Range("D34").Style = $style_D34
$ws.Range("E34").Value = '  +0.33%  '

$ws.Range("E35").Value = '  -3.51%  '

$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").Value = "'7.87"
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = '  +0.66%  '

$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("E38").Value = '  -1.61%  '

$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").Value = "'5.62"
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = '  -0.11%  '

$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").Value = "'174.27"
$ws.Range("D40").Style = $style_D40
$ws.Range("E40").Value = '  +0.69%  '

$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.0850"
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = '  -0.61%  '

$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").Value = "'5.20"
$ws.Range("D42").Style = $style_D42
$ws.Range("E42").Value = '  +0.51%  '

$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.889"
$ws.Range("D43").Style = $style_D43
$ws.Range("E43").Value = '  -0.47%  '

$ws.Range("E44").Value = '  +2.78%  '

$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").Value = "'46.11"
$ws.Range("D45").Style = $style_D45
$ws.Range("E45").Value = '  +1.41%  '

$style_D46 = $ws.Range("D46").Style
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = $style_D46
$ws.Range("E46").Value = '  +0.03%  '

$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").Value = "'2.54"
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = '  +5.48%  '

$ws.Range("E48").Value = '  -3.26%  '

$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").Value = "'25.12"
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = '  -3.48%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").Value = "'23.60"
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = '  +4.83%  '

$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").Value = "'7.17"
$ws.Range("D51").Style = $style_D51
$ws.Range("E51").Value = '  +0.47%  '
